$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in today's attendance row (row 14) for the Scrum meeting attendance roll.
$ws.Range("B14").Value = "9/19 / 1:00"

$ws.Range("C14").Value = "Google Hangout"
$ws.Range("C14").Font.Bold = $true
$ws.Range("C14").Borders.Item(8).LineStyle = -4142

$ws.Range("D14").Value = "A"
$ws.Range("E14").Value = "A"
$ws.Range("F14").Value = "A"
$ws.Range("G14").Value = "A"
$ws.Range("H14").Value = "U"
$ws.Range("I14").Value = "A"

$ws.Range("B15").Select()
